# Insert a new weekly price record row at worksheet row 245 (pushing the
# existing rows 245-354 down to 246-355), then populate the new row with
# the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 245; all subsequent rows shift down by one.
$ws.Rows("245:245").Insert()

# Populate the newly inserted row 245 with the new record's data.
$ws.Range("A245").Value = 7
$ws.Range("B245").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C245").Value = "Ñuble"
$ws.Range("D245").Value = 45134
$ws.Range("E245").Value = 16
$ws.Range("F245").Value = "Fruta"
$ws.Range("G245").Value = 100108
$ws.Range("H245").Value = "Tropicales y subtropicales"
$ws.Range("I245").Value = 100108005
$ws.Range("J245").Value = "Piña"
$ws.Range("K245").Value = "Caramelo"
$ws.Range("L245").Value = "Segunda"
$ws.Range("M245").Value = 60
$ws.Range("N245").Value = 19000
$ws.Range("O245").Value = 19000
$ws.Range("P245").Value = 19000
$ws.Range("Q245").Value = '$/caja 14 unidades'
$ws.Range("R245").Value = "Ecuador"
$ws.Range("S245").Value = 1357
$ws.Range("T245").Value = 14
